$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Header row (row 1) text updates ---
# the header row keeps the same bold/filled style (s=1); only the
# text shown in each column is replaced with the new header names
$ws.Range("A1").Value = "Tipo de Respuesta"
$ws.Range("B1").Value = "Tipo de producto"
$ws.Range("C1").Value = "Código"
$ws.Range("D1").Value = "Departamento"
$ws.Range("E1").Value = "Municipio"
$ws.Range("F1").Value = "Fecha de Elaboración"
$ws.Range("G1").Value = "Año"
$ws.Range("H1").Value = "Enlace"
$ws.Range("H1").Style = "Normal"

# --- Step 2: capture the existing link values before touching them ---
$ernUrl = $ws.Range("G2").Value2
$fceUrl = $ws.Range("H2").Value2

# remove the two existing hyperlinks (both will be recreated further right)
$ws.Range("G2").Hyperlinks.Delete()

# G2 disappears completely in the new layout
$ws.Range("G2").Clear()
# H3 disappears completely in the new layout
$ws.Range("H3").Clear()

# --- Step 3: place the urls onto their new homes: H2 (ERN) / I2 (FCE) ---
$ws.Range("H2").Value = $ernUrl
$ws.Range("I2").Value = $fceUrl

# recreate the hyperlinks (fce/I2 first so it keeps matching rId1, ern/H2 second -> rId2,
# mirroring the relationship ids used before the edit)
$ws.Hyperlinks.Add($ws.Range("I2"), $fceUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("H2"), $ernUrl) | Out-Null

# restore the same visual format used by the rest of the column (s=9) since
# Hyperlinks.Add() forces its own formatting
$ws.Range("G3").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Step 4: the note cell moves from I2 to J2 with default style and new text ---
$ws.Range("J2").Value = "por favor respetar este formato"

# --- Step 5: update the selection to mimic the target sheetView state ---
$ws.Range("A1:G1048576").Select()
